$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.710.98"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.917.98"
$ws.Range("E3").Value = "  +1.62%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.73"
$ws.Range("E5").Value = "  -1.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4941"
$ws.Range("E7").Value = "  +0.63%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3008"
$ws.Range("E8").Value = "  +2.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06781"
$ws.Range("E9").Value = "  +0.13%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.905.58"
$ws.Range("E10").Value = "  +0.93%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.23"
$ws.Range("E11").Value = "  -0.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07319"
$ws.Range("E12").Value = "  +1.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.232"
$ws.Range("E13").Value = "  +3.73%  "

# Row 14
$ws.Range("E14").Value = "  -2.85%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6754"
$ws.Range("E15").Value = "  -0.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.683.07"
$ws.Range("E16").Value = "  +0.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007984"
$ws.Range("E17").Value = "  +0.26%  "

# Row 18
$ws.Range("E18").Value = "  +3.26%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.165.19"
$ws.Range("E20").Value = "  +1.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.377"
$ws.Range("E21").Value = "  +11.50%  "

# Row 22
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "196.16"
$ws.Range("E23").Value = "  +1.96%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.342"
$ws.Range("E24").Value = "  +4.52%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.682"
$ws.Range("E25").Value = "  +3.86%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.94"
$ws.Range("E26").Value = "  +4.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.65"
$ws.Range("E27").Value = "  -2.58%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.963"
$ws.Range("E28").Value = "  +3.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.465"
$ws.Range("E29").Value = "  +4.64%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.357"
$ws.Range("E30").Value = "  +0.80%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09149"
$ws.Range("E31").Value = "  +1.28%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.080"
$ws.Range("E32").Value = "  +1.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05297"
$ws.Range("E33").Value = "  +1.91%  "

# Row 34
$ws.Range("E34").Value = "  -1.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.122"
$ws.Range("E35").Value = "  +1.08%  "

# Row 36
$ws.Range("E36").Value = "  -1.70%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01855"
$ws.Range("E37").Value = "  +1.20%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.724"
$ws.Range("E38").Value = "  +2.15%  "

# Row 39
$ws.Range("E39").Value = "  -0.31%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.087"
$ws.Range("E40").Value = "  -2.56%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4505"
$ws.Range("E41").Value = "  +2.07%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.12"
$ws.Range("E42").Value = "  +24.78%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.940"
$ws.Range("E43").Value = "  +3.62%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.00"
$ws.Range("E44").Value = "  +1.82%  "

# Row 45
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1405"
$ws.Range("E45").Value = "  +4.65%  "

# Row 46
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.18%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.723"
$ws.Range("E47").Value = "  +1.75%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.47"
$ws.Range("E48").Value = "  +5.80%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.130"
$ws.Range("E49").Value = "  +4.83%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05900"
$ws.Range("E50").Value = "  +0.77%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4046"
$ws.Range("E51").Value = "  +3.28%  "
